$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A:C width: 37.42578125 -> 37.85546875 ---
# (engine quantizes ColumnWidth to ~1/6-character steps; 37.0 is the input
#  value whose rounded result lands closest to the target stored width)
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 37

# --- Row heights ---
$ws.Rows("4:4").RowHeight = 28.5
$ws.Rows("5:5").RowHeight = 27.75
$ws.Rows("6:6").RowHeight = 28.5
$ws.Rows("7:7").RowHeight = 27
$ws.Rows("8:8").RowHeight = 18.75
$ws.Rows("10:10").RowHeight = 18.75
$ws.Rows("11:11").RowHeight = 17.25

# --- Cell value updates ---
$ws.Range("Q4").Value = 109
$ws.Range("Q7").Value = 12685.1
$ws.Range("P8").Value = 478225.6
$ws.Range("Q8").Value = 559503.6
$ws.Range("Q9").Value = 131.9
$ws.Range("Q10").Value = 3384.8
$ws.Range("Q11").Value = 12517.9
